$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text representation instead of
# being auto-converted to numbers by Excel (values like "0.7179" or "241.68").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.784.86'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '1.872.47'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = '0.7179'
$ws.Range("E5").Value = '  -2.78%  '
$ws.Range("D6").Value = '241.68'
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").Value = '0.3147'
$ws.Range("E8").Value = '  -0.33%  '
$ws.Range("D9").Value = '0.07528'
$ws.Range("E9").Value = '  +4.73%  '
$ws.Range("D10").Value = '24.51'
$ws.Range("E10").Value = '  -0.85%  '
$ws.Range("D11").Value = '0.08188'
$ws.Range("E11").Value = '  -2.07%  '
$ws.Range("D12").Value = '0.7429'
$ws.Range("E12").Value = '  -0.99%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '5.326'
$ws.Range("E13").Value = '  -1.76%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.855.57'
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").Value = '92.41'
$ws.Range("E15").Value = '  -0.17%  '
$ws.Range("D16").Value = '29.846.92'
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("D17").Value = '6.008'
$ws.Range("E17").Value = '  -0.99%  '
$ws.Range("D18").Value = '246.18'
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").Value = '0.000007918'
$ws.Range("E19").Value = '  +1.06%  '
$ws.Range("E20").Value = '  -0.80%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.145.37'
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("D24").Value = '7.713'
$ws.Range("E24").Value = '  -3.45%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '9.177'
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '0.1494'
$ws.Range("E26").Value = '  -3.69%  '
$ws.Range("D27").Value = '163.64'
$ws.Range("E27").Value = '  -0.90%  '
$ws.Range("D28").Value = '18.54'
$ws.Range("E28").Value = '  -0.62%  '
$ws.Range("D29").Value = '2.001'
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("D30").Value = '1.427'
$ws.Range("E30").Value = '  -5.34%  '
$ws.Range("D31").Value = '4.541'
$ws.Range("E31").Value = '  -0.95%  '
$ws.Range("D32").Value = '1.525'
$ws.Range("E32").Value = '  -0.75%  '
$ws.Range("D33").Value = '4.173'
$ws.Range("E33").Value = '  -2.27%  '
$ws.Range("D34").Value = '0.05451'
$ws.Range("E34").Value = '  +2.53%  '
$ws.Range("D35").Value = '1.222'
$ws.Range("E35").Value = '  -1.22%  '
$ws.Range("D36").Value = '0.7365'
$ws.Range("D37").Value = '0.9991'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").Value = '2.704'
$ws.Range("E38").Value = '  +0.24%  '
$ws.Range("D39").Value = '0.01908'
$ws.Range("E39").Value = '  -2.52%  '
$ws.Range("D40").Value = '2.733'
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("D41").Value = '0.4448'
$ws.Range("E41").Value = '  -1.46%  '
$ws.Range("D42").Value = '0.8824'
$ws.Range("E42").Value = '  +3.29%  '
$ws.Range("D43").Value = '5.991'
$ws.Range("E43").Value = '  -0.97%  '
$ws.Range("D44").Value = '71.33'
$ws.Range("E44").Value = '  -1.44%  '
$ws.Range("D45").Value = '1.003'
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '103.60'
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '1.035.90'
$ws.Range("E47").Value = '  -6.94%  '
$ws.Range("D48").Value = '7.455'
$ws.Range("E48").Value = '  -2.43%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.561'
$ws.Range("E49").Value = '  +0.69%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '1.790'
$ws.Range("E50").Value = '  -2.72%  '
$ws.Range("D51").Value = '2.035.47'
$ws.Range("E51").Value = '  +0.52%  '
